# fix figures in section 31
#
# 1) Reposition / resize the figure picture on slide 1 and lock its aspect
#    ratio (picLocks noChangeAspect="1").
# 2) Refresh the cached "datetimeFigureOut" footer-date text (8/25/20 ->
#    11/17/2020) on the slide master and on every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Picture 3 on slide 1: new position/size + lock aspect ratio.
# ---------------------------------------------------------------------
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

# Values below are expressed in points (EMU / 12700) and nudged to the
# exact double that this host's float32-backed Left/Top/Width/Height
# setters round-trip to the target EMU values:
#   x  = 365656 EMU, y  = 1099457 EMU
#   cx = 10978697 EMU, cy = 4463143 EMU
$sh.Left   = 28.791812896728516
$sh.Top    = 86.57141876220703
$sh.Width  = 864.46435546875
$sh.Height = 351.4285888671875

$sh.LockAspectRatio = -1   # msoTrue -> <a:picLocks noChangeAspect="1"/>

# ---------------------------------------------------------------------
# 2) Update the cached date placeholder text everywhere it appears.
# ---------------------------------------------------------------------
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Name -like "Date Placeholder*") {
            if ($candidate.HasTextFrame) {
                $candidate.TextFrame.TextRange.Text = "11/17/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}
